# Regenerate save_data column G ("K" - strikeouts) values to replace the
# previous "Strike#" derived numbers with the correct K counts.
# (commit: "regen save_data to use K instead of Strike#, regen std/mean,
#  calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> new K value, per the target diff (column G, rows 2-39).
$kValues = [ordered]@{
    2  = 4
    3  = 8
    4  = 7
    5  = 8
    6  = 8
    7  = 2
    8  = 6
    9  = 9
    10 = 3
    11 = 8
    12 = 6
    13 = 7
    14 = 9
    15 = 11
    16 = 5
    17 = 3
    18 = 6
    19 = 7
    20 = 3
    21 = 6
    22 = 9
    23 = 6
    24 = 10
    25 = 8
    26 = 6
    27 = 9
    28 = 5
    29 = 4
    30 = 8
    31 = 9
    32 = 5
    33 = 6
    34 = 6
    35 = 6
    36 = 2
    37 = 4
    38 = 4
    39 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
